$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Sprint backlog"

# --- Row 20: copy formatting from row 19 (plain bordered style) ---
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Description of tasks per week"
$ws.Range("C20").Value = "Darius"
$ws.Range("D20").Value = $null
$ws.Range("E20").Value = $null
$ws.Range("F20").Value = $null
$ws.Range("G20").Value = 1

# --- Row 21: copy formatting from row 7 (shaded "Presentation" row) ---
$ws.Range("A7:G7").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)
$ws.Range("A21").Value = 11
$ws.Range("C21").Value = "Hayden"
$ws.Range("D21").Value = $null
$ws.Range("E21").Value = $null
$ws.Range("F21").Value = $null
$ws.Range("G21").Value = 1

# --- Row 22: same shaded style as row 21 ---
$ws.Range("A7:G7").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$ws.Range("A22").Value = 12
$ws.Range("B22").Value = "Group Evaluation"
$ws.Range("C22").Value = "Darius"
$ws.Range("D22").Value = $null
$ws.Range("E22").Value = $null
$ws.Range("F22").Value = $null
$ws.Range("G22").Value = 1

# The task description for row 21 is filled in last (matches shared-string order)
$ws.Range("B21").Value = "Sprint Chart Generator"

# Column widths for D and E changed from a shared width to two distinct widths
$ws.Range("D1").ColumnWidth = 10.25
$ws.Range("E1").ColumnWidth = 9.15

# Update the active cell selection to match the saved state
$ws.Range("C26").Select()
